$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

# "Ngày tiêm" moves from I1 to the new last column J1,
# and I1 becomes the newly introduced "Ngày sinh" column.
$ws.Cells.Item(1, 10).Value = "Ngày tiêm"
$ws.Cells.Item(1, 9).Value = "Ngày sinh"

$ws.Range("G4").Select()
